{"js": "// Office.js (Word JavaScript API) edit script.\n// Implements: \"Added file counts to the submission instructions\"\n//   1. \"Upload the following to the ...\" -> \"Upload the following 7 files to the ...\"\n//   2. \"A zip file containing the two files (.html and .js) for part 1.\"\n//        -> \"Two files (.html and .js) for part 1.\"\n//   3. \"A zip file containing the four files for part 2.\"\n//        -> \"Four files (2 html and 2 js) for part 2.\"\n//   4. Moves the \"_GoBack\" bookmark so it spans from the start of the\n//      \"Upload the following...\" paragraph through the end of the\n//      \"...filled in by you.\" paragraph (instead of sitting mid-word\n//      inside \"with\").\n\nconst body = context.document.body;\n\n// --- Text edits -----------------------------------------------------\n// Merge \"...your lab partner wi\" + (bookmark) + \"th the ...\" into a single\n// run so the stray \"_GoBack\" bookmark collapses to one spot (right after\n// \"with the \") before it gets deleted/recreated below.\nconst withTheHits = body.search(\" your lab partner with the \", { matchCase: true });\nwithTheHits.load(\"items\");\nawait context.sync();\nif (withTheHits.items.length > 0) {\n  withTheHits.items[0].insertText(\" your lab partner with the \", \"Replace\");\n}\nawait context.sync();\n\n// Paragraph: \"Upload the following to the Lab Production Version assignment:\"\nconst uploadHits = body.search(\"Upload the following to the \", { matchCase: true });\nuploadHits.load(\"items\");\nawait context.sync();\nif (uploadHits.items.length > 0) {\n  uploadHits.items[0].insertText(\"Upload the following 7 files to the \", \"Replace\");\n}\n\n// Paragraph: \"A zip file containing the two files (.html and .js) for part 1.\"\nconst part1Hits = body.search(\"A zip file containing the two files (.html and .\", { matchCase: true });\npart1Hits.load(\"items\");\nawait context.sync();\nif (part1Hits.items.length > 0) {\n  part1Hits.items[0].insertText(\"Two files (.html and .\", \"Replace\");\n}\n\n// Paragraph: \"A zip file containing the four files for part 2.\"\nconst part2Hits = body.search(\"A zip file containing the four files for part 2.\", { matchCase: true });\npart2Hits.load(\"items\");\nawait context.sync();\nif (part2Hits.items.length > 0) {\n  part2Hits.items[0].insertText(\"Four files (2 html and 2 js) for part 2.\", \"Replace\");\n}\n\nawait context.sync();\n\n// --- Move the \"_GoBack\" bookmark ------------------------------------\n// Remove the old bookmark (currently sitting mid-word inside \"with\").\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// Re-resolve the paragraphs after the text edits above, then re-create the\n// bookmark spanning from the start of the \"Upload the following...\"\n// paragraph to the end of the \"...filled in by you.\" paragraph.\nconst uploadPara = body.search(\"Upload the following 7 files to the\", { matchCase: true });\nuploadPara.load(\"items\");\nconst partnerPara = body.search(\"filled in by you.\", { matchCase: true });\npartnerPara.load(\"items\");\nawait context.sync();\n\nif (uploadPara.items.length > 0 && partnerPara.items.length > 0) {\n  const startRange = uploadPara.items[0].paragraphs.getFirst().getRange(\"Start\");\n  const endRange = partnerPara.items[0].paragraphs.getLast().getRange(\"End\");\n  const fullRange = startRange.expandTo(endRange);\n  fullRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Implements: \"Added file counts to the submission instructions\"\n#   1. \"Upload the following to the ...\" -> \"Upload the following 7 files to the ...\"\n#   2. \"A zip file containing the two files (.html and .js) for part 1.\"\n#        -> \"Two files (.html and .js) for part 1.\"\n#   3. \"A zip file containing the four files for part 2.\"\n#        -> \"Four files (2 html and 2 js) for part 2.\"\n#   4. Moves the \"_GoBack\" bookmark so it spans from the start of the\n#      \"Upload the following...\" paragraph through the end of the\n#      \"...filled in by you.\" paragraph (instead of sitting mid-word\n#      inside \"with\").\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n\n# --- Remove the old \"_GoBack\" bookmark first (it currently sits mid-word\n# inside \"with\", between \"wi\" and \"th the\") so it doesn't end up orphaned\n# or duplicated once the text around it is edited.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Text edits -------------------------------------------------------\n# Re-merge \" your lab partner wi\" + \"th the \" (split where the old\n# bookmark used to sit) back into a single run of identical text.\nReplace-Text \" your lab partner with the \" \" your lab partner with the \"\nReplace-Text \"Upload the following to the \" \"Upload the following 7 files to the \"\nReplace-Text \"A zip file containing the two files (.html and .\" \"Two files (.html and .\"\nReplace-Text \"A zip file containing the four files for part 2.\" \"Four files (2 html and 2 js) for part 2.\"\n\n# --- Re-create the \"_GoBack\" bookmark spanning from the start of the\n# \"Upload the following...\" paragraph to the end of the\n# \"...filled in by you.\" paragraph.\n$uploadPara = $null\n$partnerPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"Upload the following*\") {\n        $uploadPara = $p\n    }\n    if ($t -like \"*filled in by you.*\") {\n        $partnerPara = $p\n    }\n}\n\nif ($uploadPara -ne $null -and $partnerPara -ne $null) {\n    $startPos = $uploadPara.Range.Start\n    $endPos = $partnerPara.Range.End\n    $bookmarkRange = $d.Range($startPos, $endPos)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n}\n"}
